$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 6
$ws.Range("K8").Value = 18
$ws.Range("M8").Value = 121

$ws.Range("H31").Value = 277.33334
$ws.Range("I31").Value = 277.33334
$ws.Range("K31").Value = 832.0000200000001
$ws.Range("M31").Value = -602.0000200000001

$ws.Range("H39").Value = 190.90909
$ws.Range("I39").Value = 130
$ws.Range("K39").Value = 390
$ws.Range("M39").Value = -94

$ws.Range("H58").Value = 470.375
$ws.Range("I58").Value = 81.666664
$ws.Range("K58").Value = 244.999992
$ws.Range("M58").Value = -94.99999199999999

$ws.Range("H74").Value = 2328.4285
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872

$ws.Range("H76").Value = 5874.25
$ws.Range("I76").Value = 5832.6665
$ws.Range("K76").Value = 5832.6665
$ws.Range("M76").Value = -5517.6665

$ws.Range("H77").Value = 2328.4285
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360

$ws.Range("H79").Value = 5874.25
$ws.Range("I79").Value = 5832.6665
$ws.Range("K79").Value = 5832.6665
$ws.Range("M79").Value = -4740.6665

$ws.Range("H96").Value = 4931.25
$ws.Range("I96").Value = 4931.25
$ws.Range("K96").Value = 14793.75
$ws.Range("M96").Value = -13420.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 290.5
$ws.Range("J4").Value = 424.5
$ws.Range("L4").Value = 424.5
$ws.Range("N4").Value = -656.5

$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 200
$ws.Range("K15").Value = 200
$ws.Range("M15").Value = 150

$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60630

$ws.Range("H63").Value = 3981.3333
$ws.Range("I63").Value = 4263.143
$ws.Range("J63").Value = 2995
$ws.Range("K63").Value = 4263.143
$ws.Range("L63").Value = 2995
$ws.Range("M63").Value = -3577.143
$ws.Range("N63").Value = -4367

$ws.Range("H66").Value = 3981.3333
$ws.Range("I66").Value = 4263.143
$ws.Range("J66").Value = 2995
$ws.Range("K66").Value = 21315.715
$ws.Range("L66").Value = 14975
$ws.Range("M66").Value = -17883.715
$ws.Range("N66").Value = -21839

$ws.Range("H88").Value = 2999
$ws.Range("J88").Value = 2999.5
$ws.Range("L88").Value = 2999.5
$ws.Range("N88").Value = -3811.5

$ws.Range("H91").Value = 2999
$ws.Range("J91").Value = 2999.5
$ws.Range("L91").Value = 2999.5
$ws.Range("N91").Value = -5807.5

$ws.Range("H92").Value = 70000
$ws.Range("I92").Value = 90000
$ws.Range("J92").Value = 50000
$ws.Range("K92").Value = 90000
$ws.Range("L92").Value = 50000
$ws.Range("M92").Value = -87504
$ws.Range("N92").Value = -54992

$ws.Range("H95").Value = 17264.834
$ws.Range("J95").Value = 17264.834
$ws.Range("L95").Value = 17264.834
$ws.Range("N95").Value = -22756.834

$ws.Range("H102").Value = 1875
$ws.Range("I102").Value = 1875
$ws.Range("K102").Value = 1875
$ws.Range("M102").Value = -253

$ws.Range("H106").Value = 15370
$ws.Range("J106").Value = 15370
$ws.Range("L106").Value = 15370
$ws.Range("N106").Value = -17894

$ws.Range("H112").Value = 31193
$ws.Range("J112").Value = 31193
$ws.Range("L112").Value = 31193
$ws.Range("N112").Value = -34147

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1507
$ws.Range("J86").Value = 1507
$ws.Range("L86").Value = 1507
$ws.Range("N86").Value = -3753

$ws.Range("H89").Value = 1507
$ws.Range("J89").Value = 1507
$ws.Range("L89").Value = 7535
$ws.Range("N89").Value = -18767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 296.7143
$ws.Range("I107").Value = 279.5
$ws.Range("K107").Value = 279.5
$ws.Range("M107").Value = 1640.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 725.76
$ws.Range("I2").Value = 336.66666
$ws.Range("J2").Value = 2768.5
$ws.Range("K2").Value = 2019.99996
$ws.Range("L2").Value = 16611
$ws.Range("M2").Value = -1906.99996
$ws.Range("N2").Value = -16837

$ws.Range("H104").Value = 499.5
$ws.Range("J104").Value = 499.33334
$ws.Range("L104").Value = 1498.00002
$ws.Range("N104").Value = -6740.000019999999

$ws.Range("H131").Value = 3361.3333
$ws.Range("I131").Value = 4222.4
$ws.Range("J131").Value = 2285
$ws.Range("K131").Value = 12667.2
$ws.Range("L131").Value = 6855
$ws.Range("M131").Value = -7627.199999999999
$ws.Range("N131").Value = -16935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 3013.75
$ws.Range("I57").Value = 3013.75
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 3013.75
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -2193.75
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 6125.25
$ws.Range("I80").Value = 4001
$ws.Range("K80").Value = 4001
$ws.Range("M80").Value = -3003

$ws.Range("H83").Value = 6125.25
$ws.Range("I83").Value = 4001
$ws.Range("K83").Value = 20005
$ws.Range("M83").Value = -15013

$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5633.3335
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376

$ws.Range("H55").Value = 2779.375
$ws.Range("J55").Value = 2947
$ws.Range("L55").Value = 2947
$ws.Range("N55").Value = -3293

$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344

$ws.Range("H110").Value = 37995.668
$ws.Range("J110").Value = 37995.668
$ws.Range("L110").Value = 37995.668
$ws.Range("N110").Value = -46175.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

$ws.Range("H100").Value = 1110.125
$ws.Range("I100").Value = 628.8
$ws.Range("K100").Value = 1257.6
$ws.Range("M100").Value = -716.5999999999999
